$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 779.9074730793836
$ws.Range("C2").Value = 1930.30262026198
$ws.Range("D2").Value = 3570.578798461907
